# Auto-generated edit applying numeric market-data refresh to Sheets/Midgardsormr_Profits.xlsx
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 823.71185
$ws.Range("J17").Value = 828.0702
$ws.Range("L17").Value = 2484.2106
$ws.Range("N17").Value = -2820.2106
$ws.Range("H33").Value = 404.33334
$ws.Range("I33").Value = 109.23529
$ws.Range("J33").Value = 1121
$ws.Range("K33").Value = 109.23529
$ws.Range("L33").Value = 1121
$ws.Range("M33").Value = 119.76471
$ws.Range("N33").Value = -1579
$ws.Range("H86").Value = 1995.8
$ws.Range("I86").Value = 1672.0714
$ws.Range("J86").Value = 2751.1667
$ws.Range("K86").Value = 1672.0714
$ws.Range("L86").Value = 2751.1667
$ws.Range("M86").Value = -549.0714
$ws.Range("N86").Value = -4997.1667
$ws.Range("H89").Value = 1995.8
$ws.Range("I89").Value = 1672.0714
$ws.Range("J89").Value = 2751.1667
$ws.Range("K89").Value = 8360.357
$ws.Range("L89").Value = 13755.8335
$ws.Range("M89").Value = -2744.357
$ws.Range("N89").Value = -24987.8335
$ws.Range("H106").Value = 3264.3635
$ws.Range("I106").Value = 3264.3635
$ws.Range("K106").Value = 3264.3635
$ws.Range("M106").Value = -2633.3635
$ws.Range("H112").Value = 4837.394
$ws.Range("I112").Value = 1416
$ws.Range("J112").Value = 5179.533
$ws.Range("K112").Value = 4248
$ws.Range("L112").Value = 15538.599
$ws.Range("M112").Value = -3140
$ws.Range("N112").Value = -17754.599
$ws.Range("H137").Value = 20381.941
$ws.Range("I137").Value = 25723
$ws.Range("K137").Value = 77169
$ws.Range("M137").Value = -74619
$ws.Range("H138").Value = 36051.566
$ws.Range("I138").Value = 2093.6667
$ws.Range("J138").Value = 86988.414
$ws.Range("K138").Value = 6281.000100000001
$ws.Range("L138").Value = 260965.242
$ws.Range("M138").Value = -1141.000100000001
$ws.Range("N138").Value = -271245.242

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26589.744
$ws.Range("I32").Value = 30606.135
$ws.Range("J32").Value = 1822
$ws.Range("K32").Value = 30606.135
$ws.Range("L32").Value = 1822
$ws.Range("M32").Value = -30319.135
$ws.Range("N32").Value = -2396
$ws.Range("H80").Value = 32500.25
$ws.Range("H83").Value = 32500.25

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2138.15
$ws.Range("I94").Value = 1461
$ws.Range("K94").Value = 1461
$ws.Range("M94").Value = -1010

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2391.25
$ws.Range("I16").Value = 1600
$ws.Range("J16").Value = 2866
$ws.Range("K16").Value = 1600
$ws.Range("L16").Value = 2866
$ws.Range("M16").Value = -1313
$ws.Range("N16").Value = -3440
$ws.Range("H113").Value = 2391.25
$ws.Range("I113").Value = 1600
$ws.Range("J113").Value = 2866
$ws.Range("K113").Value = 1600
$ws.Range("L113").Value = 2866
$ws.Range("M113").Value = 570
$ws.Range("N113").Value = -7206
$ws.Range("H132").Value = 38113.89
$ws.Range("I132").Value = 46348.953
$ws.Range("J132").Value = 1879.6
$ws.Range("K132").Value = 139046.859
$ws.Range("L132").Value = 5638.799999999999
$ws.Range("M132").Value = -136516.859
$ws.Range("N132").Value = -10698.8

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2180
$ws.Range("I3").Value = 2180
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 6540
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -6428
$ws.Range("N3").ClearContents()
$ws.Range("H5").Value = 1007.58826
$ws.Range("I5").Value = 923.5
$ws.Range("J5").Value = 1127.7142
$ws.Range("K5").Value = 2770.5
$ws.Range("L5").Value = 3383.1426
$ws.Range("M5").Value = -2658.5
$ws.Range("N5").Value = -3607.1426
$ws.Range("H113").Value = 798.7692
$ws.Range("I113").Value = 342.5
$ws.Range("J113").Value = 881.7273
$ws.Range("K113").Value = 1027.5
$ws.Range("L113").Value = 2645.1819
$ws.Range("M113").Value = 1142.5
$ws.Range("N113").Value = -6985.1819
$ws.Range("H124").Value = 4456.522
$ws.Range("I124").Value = 833.3333
$ws.Range("J124").Value = 5000
$ws.Range("K124").Value = 2499.9999
$ws.Range("L124").Value = 15000
$ws.Range("M124").Value = 2410.0001
$ws.Range("N124").Value = -24820
$ws.Range("H125").Value = 800
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H129").Value = 2521.4736
$ws.Range("I129").Value = 1690.4615
$ws.Range("J129").Value = 4322
$ws.Range("K129").Value = 5071.3845
$ws.Range("L129").Value = 12966
$ws.Range("M129").Value = -71.38450000000012
$ws.Range("N129").Value = -22966
$ws.Range("H133").Value = 2825.2
$ws.Range("I133").Value = 2741.2856
$ws.Range("K133").Value = 8223.856800000001
$ws.Range("M133").Value = -3163.856800000001
$ws.Range("H135").Value = 1007.58826
$ws.Range("I135").Value = 923.5
$ws.Range("J135").Value = 1127.7142
$ws.Range("K135").Value = 8311.5
$ws.Range("L135").Value = 10149.4278
$ws.Range("M135").Value = -5776.5
$ws.Range("N135").Value = -15219.4278
$ws.Range("H137").Value = 2987.389
$ws.Range("I137").Value = 1982.2858
$ws.Range("J137").Value = 3627
$ws.Range("K137").Value = 5946.857400000001
$ws.Range("L137").Value = 10881
$ws.Range("M137").Value = -846.8574000000008
$ws.Range("N137").Value = -21081

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 174.90909
$ws.Range("I2").Value = 138.27272
$ws.Range("K2").Value = 138.27272
$ws.Range("M2").Value = -25.27271999999999
$ws.Range("H26").Value = 38000
$ws.Range("J26").Value = 38000
$ws.Range("L26").Value = 38000
$ws.Range("N26").Value = -38560
$ws.Range("H50").Value = 38000
$ws.Range("J50").Value = 38000
$ws.Range("L50").Value = 38000
$ws.Range("N50").Value = -38996
$ws.Range("H122").Value = 3348.4546
$ws.Range("I122").Value = 3225.389
$ws.Range("K122").Value = 9676.167000000001
$ws.Range("M122").Value = -7226.167000000001
$ws.Range("H126").Value = 2939.1333
$ws.Range("I126").Value = 2576.6667
$ws.Range("J126").Value = 3482.8333
$ws.Range("K126").Value = 7730.000100000001
$ws.Range("L126").Value = 10448.4999
$ws.Range("M126").Value = -5260.000100000001
$ws.Range("N126").Value = -15388.4999
$ws.Range("H132").Value = 2625.7727
$ws.Range("I132").Value = 2512.7144
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 7538.1432
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -5008.1432
$ws.Range("N132").Value = -20060

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1000.2308
$ws.Range("I16").Value = 1437.5
$ws.Range("K16").Value = 1437.5
$ws.Range("M16").Value = -1267.5
$ws.Range("H61").Value = 1156.4
$ws.Range("I61").Value = 1169.25
$ws.Range("J61").Value = 1105
$ws.Range("K61").Value = 1169.25
$ws.Range("L61").Value = 1105
$ws.Range("M61").Value = -967.25
$ws.Range("N61").Value = -1509
$ws.Range("H82").Value = 1474.6428
$ws.Range("I82").Value = 1356.4286
$ws.Range("J82").Value = 1592.8572
$ws.Range("K82").Value = 1356.4286
$ws.Range("L82").Value = 1592.8572
$ws.Range("M82").Value = -995.4286
$ws.Range("N82").Value = -2314.8572
$ws.Range("H85").Value = 1474.6428
$ws.Range("I85").Value = 1356.4286
$ws.Range("J85").Value = 1592.8572
$ws.Range("K85").Value = 1356.4286
$ws.Range("L85").Value = 1592.8572
$ws.Range("M85").Value = -108.4286
$ws.Range("N85").Value = -4088.8572
$ws.Range("H93").Value = 1824.3077
$ws.Range("I93").Value = 1535.3334
$ws.Range("J93").Value = 2474.5
$ws.Range("K93").Value = 1535.3334
$ws.Range("L93").Value = 2474.5
$ws.Range("M93").Value = -287.3334
$ws.Range("N93").Value = -4970.5
$ws.Range("H113").Value = 1156.4
$ws.Range("I113").Value = 1169.25
$ws.Range("J113").Value = 1105
$ws.Range("K113").Value = 1169.25
$ws.Range("L113").Value = 1105
$ws.Range("M113").Value = 1000.75
$ws.Range("N113").Value = -5445
$ws.Range("H122").Value = 3882.2122
$ws.Range("I122").Value = 3265.625
$ws.Range("J122").Value = 5526.4443
$ws.Range("K122").Value = 9796.875
$ws.Range("L122").Value = 16579.3329
$ws.Range("M122").Value = -7346.875
$ws.Range("N122").Value = -21479.3329
$ws.Range("H132").Value = 1946.9474
$ws.Range("I132").Value = 908.0909
$ws.Range("K132").Value = 2724.2727
$ws.Range("M132").Value = -194.2727

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1727.1428
$ws.Range("J100").Value = 1592.6666
$ws.Range("L100").Value = 3185.3332
$ws.Range("N100").Value = -4267.3332
$ws.Range("H122").Value = 194291.42
$ws.Range("I122").Value = 329928.28
$ws.Range("K122").Value = 989784.8400000001
$ws.Range("M122").Value = -987334.8400000001
$ws.Range("H126").Value = 253889.16
$ws.Range("I126").Value = 3471.8667
$ws.Range("J126").Value = 1005141
$ws.Range("K126").Value = 10415.6001
$ws.Range("L126").Value = 3015423
$ws.Range("M126").Value = -7945.6001
$ws.Range("N126").Value = -3020363
